$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update A3 from 1 to 6 (rest of row 3 stays the same)
$ws.Range("A3").Value = 6

# Row 4 (new)
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "x"
$ws.Range("D4").Value = "x"
$ws.Range("F4").Value = "x"
$ws.Range("G4").Value = "x"
$ws.Range("I4").Value = "x"
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = 6
$ws.Range("P4").Value = 2
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = "x"

# Row 5 (new)
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = 2
$ws.Range("D5").Value = "x"
$ws.Range("F5").Value = "x"
$ws.Range("H5").Value = "x"
$ws.Range("J5").Value = "x"
$ws.Range("K5").Value = 12
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 4
$ws.Range("U5").Value = "x"

# Update the active cell / selection to C7
$ws.Range("C7").Select()
